# Append/merge latest scrape results into the "ランサーズ" sheet.
# The sheet keeps listings sorted descending by column G (優先度スコア).
# This run found 3 new listings (inserted so the sheet stays sorted) and
# refreshed the "取得日時" timestamp on every existing listing. One
# previously-seen listing's price/terms text ("取引期間") was updated too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-11-20 01:18:39"

# Drop all existing hyperlinks up front - row inserts below do not shift the
# <hyperlinks> anchors automatically, so it is simplest to rebuild them once
# every row is in its final place.
$ws.Hyperlinks.Delete()

# Insert the 3 new rows top-to-bottom so each subsequent row index below is
# still relative to the sheet *after* the previous inserts.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(11).Insert()

# New row 2: 【AI開発】生成AI・RAGシステム構築パートナー募集
$ws.Cells.Item(2, 1).Value = $newTimestamp
$ws.Cells.Item(2, 2).Value = "【AI開発】生成AI・RAGシステム構築パートナー募集"
$ws.Cells.Item(2, 3).Value = "システム開発"
$ws.Cells.Item(2, 4).Value = "1,000 ~ 5,000 円 / 固定"
$ws.Cells.Item(2, 5).Value = "期限情報なし"
$ws.Cells.Item(2, 6).Value = "https://www.lancers.jp/work/detail/5437447"
$ws.Cells.Item(2, 7).Value = 375
$ws.Cells.Item(2, 8).Value = "🔥AI,Ai ◆開発"

# New row 8: 【急募】価格更新サイトにエクセルアップロード後、内容を更新するプログラム作成依頼
$ws.Cells.Item(8, 1).Value = $newTimestamp
$ws.Cells.Item(8, 2).Value = "【急募】価格更新サイトにエクセルアップロード後、内容を更新するプログラム作成依頼"
$ws.Cells.Item(8, 3).Value = "システム開発"
$ws.Cells.Item(8, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(8, 5).Value = "期限情報なし"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5437655"
$ws.Cells.Item(8, 7).Value = 38
$ws.Cells.Item(8, 8).Value = "◇サイト"

# New row 11: 【SESエンジニア募集】多様なプロジェクトに参画可能!
$ws.Cells.Item(11, 1).Value = $newTimestamp
$ws.Cells.Item(11, 2).Value = "【SESエンジニア募集】多様なプロジェクトに参画可能!"
$ws.Cells.Item(11, 3).Value = "システム開発"
$ws.Cells.Item(11, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(11, 5).Value = "期限情報なし"
$ws.Cells.Item(11, 6).Value = "https://www.lancers.jp/work/detail/5437544"
$ws.Cells.Item(11, 7).Value = 25

# Refresh the "取得日時" timestamp on every pre-existing listing row (the
# new rows above already carry it too).
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Existing listing "【急募】Android用のライブ壁紙アプリ開発エンジニアを探しています!"
# (now row 6) picked up one more day of transaction time on this crawl.
$ws.Cells.Item(6, 4).Value = "100,000 円 ~ 200,000 円 / 募集期間 7 日、取引期間 1 日"

# Rebuild the hyperlinks for every row's URL cell (column F), matching the
# original Hyperlink cell style.
for ($r = 2; $r -le 15; $r++) {
    $target = $ws.Cells.Item($r, 6).Value2
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $target)
}
$ws.Range("F2:F15").Style = "Hyperlink"
